$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2,3) {
    $ws.Range("D$r").Value = 0.09119999999999999
    $ws.Range("E$r").Value = 0.227
    $ws.Range("G$r").Value = 0.4047151277013752
    $ws.Range("H$r").Value = 0.4047151277013752
    $ws.Range("I$r").Value = 0.3025540275049116
    $ws.Range("J$r").Value = 0.2306483300589391
    $ws.Range("K$r").Value = 1.76
    $ws.Range("L$r").Value = 0.3457760314341847
    $ws.Range("M$r").Value = 0.507
    $ws.Range("N$r").Value = 0.0325
    $ws.Range("O$r").Value = 0.2880681818181818
    $ws.Range("P$r").Value = 0.507
    $ws.Range("Q$r").Value = 0.0325
    $ws.Range("R$r").Value = 0.2880681818181818
    $ws.Range("U$r").Value = 0.841
    $ws.Range("V$r").Value = 0.05391025641025641
    $ws.Range("W$r").Value = 0.2117930204572804
    $ws.Range("X$r").Value = 0.06373606207847107
    $ws.Range("Y$r").Value = 0.1480569583788093
    $ws.Range("Z$r").Value = 0.6632786030753192
    $ws.Range("AA$r").Value = 0.1529841021631483
    $ws.Range("AB$r").Value = 0.06373606207847107
    $ws.Range("AC$r").Value = 0.08924804008467721
    $ws.Range("AG$r").Value = -0.841
    $ws.Range("AJ$r").Value = -0.05698218036452334
    $ws.Range("AK$r").Value = -0.09083054325521114
    $ws.Range("AM$r").Value = -0.6919999999999999
    $ws.Range("AP$r").Value = -0.4247474747474748
    $ws.Range("AQ$r").Value = -2.225433526011561
}
